# Fruta / hortaliza, semanal
# Insert a new weekly price-observation row for "Feria Lagunitas de Puerto Montt - Zanahoria"
# right after the existing row 114 (new row 115), shifting all following rows down by one.
# The dimension grows from A1:R226 to A1:R227.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 115 -- everything currently at
# row 115 and below moves down to make room (row 226 becomes row 227).
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new weekly observation.
$ws.Range("A115").Value = 4
$ws.Range("B115").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C115").Value = "Los Lagos"
$ws.Range("D115").Value = 44484
$ws.Range("E115").Value = 10
$ws.Range("F115").Value = 100114013
$ws.Range("G115").Value = "Zanahoria"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 850
$ws.Range("K115").Value = 12000
$ws.Range("L115").Value = 12000
$ws.Range("M115").Value = 12000
$ws.Range("N115").Value = "$/saco 20 kilos"
$ws.Range("O115").Value = "Región de Ñuble"
$ws.Range("P115").Value = 600
$ws.Range("Q115").Value = 20
$ws.Range("R115").Value = "Hortaliza"
